$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the LiveSLR version/build copyright string (row 2, column B)
$ws.Range("B2").Value = "Copyright @ 2023 Cytel Inc. LiveSLR 4.0.0.0 - Build #50015"

# Move the active selection to B3 (matches the new selection in the saved file)
$ws.Range("B3").Select()

# Restore a "normal" (non-maximized) window size/position like the new workbookView
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Left = 5760
$win.Top = 3432
$win.Width = 17280
$win.Height = 9072
